$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Add the new SQI / "No gaps" row under the existing content ---
$ws.Range("A16").Value = "SQI"
$ws.Range("C16").Value = "No gaps"

# --- Switch the sheet's content over from Arial to Calibri ---
$ws.Range("A1:C2").Font.Name = "Calibri"
$ws.Range("C6").Font.Name = "Calibri"
$ws.Range("A16").Font.Name = "Calibri"
$ws.Range("C16").Font.Name = "Calibri"

# --- Match the saved selection from the author's last edit ---
$ws.Range("C12").Select() | Out-Null
